$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Cole Anthony", "PG", "Orlando Magic"),
    @("Tristan da Silva", "SF", "Orlando Magic"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
